$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# ---------------------------------------------------------------------------
# Add new rows 15-29: UI variable mapping for IAM arrays in the Physical
# Trough model (csp_dtr_sca_iam* -> IAMs_n[i], and IamF0/1/2 -> IAM_matrix)
# ---------------------------------------------------------------------------

$rows = @(
    @("Changed name","number","csp_dtr_sca_iam0_1","IAMs_1[0]","Physical Trough Collector Type 1","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_1","IAMs_1[1]","Physical Trough Collector Type 1","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_1","IAMs_1[2]","Physical Trough Collector Type 1","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam0_2","IAMs_2[0]","Physical Trough Collector Type 2","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_2","IAMs_2[1]","Physical Trough Collector Type 2","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_2","IAMs_2[2]","Physical Trough Collector Type 2","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam0_3","IAMs_3[0]","Physical Trough Collector Type 3","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_3","IAMs_3[1]","Physical Trough Collector Type 3","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_3","IAMs_3[2]","Physical Trough Collector Type 3","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam0_4","IAMs_4[0]","Physical Trough Collector Type 4","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_4","IAMs_4[1]","Physical Trough Collector Type 4","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_4","IAMs_4[2]","Physical Trough Collector Type 4","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","array","IamF0","IAM_matrix","Physical Trough Collector Header","combining collector IAM coef. Arrays into 1 output matrix","Ty"),
    @("Changed name","array","IamF1","IAM_matrix","Physical Trough Collector Header","combining collector IAM coef. Arrays into 1 output matrix","Ty"),
    @("Changed name","array","IamF2","IAM_matrix","Physical Trough Collector Header","combining collector IAM coef. Arrays into 1 output matrix","Ty")
)

$startRow = 15
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}

# ---------------------------------------------------------------------------
# Extend the "Types" list data validation on column A down to row 58
# ---------------------------------------------------------------------------
$ws.Range("A2:A50").Validation.Delete()
$ws.Range("A2:A58").Validation.Add(3, 1, 1, "Types")
$ws.Range("A2:A58").Validation.IgnoreBlank = $false

# ---------------------------------------------------------------------------
# Widen column F to fit the new, longer default-value/reason text
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 50.14

# ---------------------------------------------------------------------------
# Move the instructional callout shape down below the newly added rows
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 435.75
$shp.Left = $shp.Left - 5.25

# ---------------------------------------------------------------------------
# Update the active selection to reflect where editing left off
# ---------------------------------------------------------------------------
$ws.Range("D29").Select()
